$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.793.62"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.635.17"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.31"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.640.87"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "1.860.20"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "0.0₃0769"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.87"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "25.805.09"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.39"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.98"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.81"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.905"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "1.131.17"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.40"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.808"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "1.769.37"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.23"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.417"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.55"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("E51").Value = "  +2.70%  "
